$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("U9")
$cmt = $c.AddComment("test")
Write-Host "comment added"
Write-Host $c.Comment.Text()
